# After setting Fp and Fl lower bound to 0, got "Better" results
# Update the 16x3 grid of result values on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @(2632.1748002241775,2133.4651017085321,2059.6828857739411),
    @(2424.3275760297829,1636.0682234783621,1818.6760318077977),
    @(2891.7242165852699,2193.7484676388208,2138.8915593313782),
    @(2715.7900851241898,2202.8507404223883,2210.8096086272531),
    @(2940.5293546419316,2148.0221079168919,2198.8077506115974),
    @(2785.7146696081741,2303.5575870482417,2387.5653411731241),
    @(2398.1173670810554,2115.3289376251591,2074.1092841307468),
    @(2882.0547411567204,2506.259365367097,2322.0381252657962),
    @(3159.4912238408479,2252.6815123713995,2035.8815171824333),
    @(2442.5676296694178,1443.4470074285916,1577.8235036845329),
    @(2285.2918750135318,1674.3145887326211,1616.9537472491388),
    @(3074.4885462675111,2504.906804466516,2362.1270015854234),
    @(3145.5104873942182,2517.3901677236399,2258.3211228889404),
    @(3159.7893414154605,2565.8316250095054,2289.3730088992443),
    @(3128.9397244237757,2515.756915603733,2389.3836037928622),
    @(3012.3636228073856,2161.4222517659405,1838.133979203845)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $rowNum = $i + 1
    $ws.Range("A$rowNum").Value = $newData[$i][0]
    $ws.Range("B$rowNum").Value = $newData[$i][1]
    $ws.Range("C$rowNum").Value = $newData[$i][2]
}
